$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.890.29"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "2.934.14"
$ws.Range("E3").Value = "  +3.43%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "352.45"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "112.09"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "39.38"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("D11").Value = "0.0890"
$ws.Range("E11").Value = "  +4.53%  "
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("D13").Value = "20.06"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "7.77"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "3.395.90"
$ws.Range("E15").Value = "  +3.45%  "
$ws.Range("D16").Value = "2.937.64"
$ws.Range("E16").Value = "  +3.47%  "
$ws.Range("D17").Value = "0.986"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "51.974.87"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  -4.63%  "
$ws.Range("D20").Value = "7.61"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "14.23"
$ws.Range("E21").Value = "  +6.52%  "
$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").Value = "71.23"
$ws.Range("E23").Value = "  +0.97%  "
$ws.Range("D24").Value = "268.59"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "2.79"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("E26").Value = "  +11.56%  "
$ws.Range("D27").Value = "26.96"
$ws.Range("E27").Value = "  +2.41%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "7.25"
$ws.Range("E29").Value = "  +14.66%  "
$ws.Range("E30").Value = "  +15.66%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("D32").Value = "6.27"
$ws.Range("E32").Value = "  +10.80%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("E38").Value = "  +3.23%  "
$ws.Range("D39").Value = "18.74"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("E41").Value = "  +6.40%  "
$ws.Range("D42").Value = "0.117"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "23.26"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("E45").Value = "  +1.57%  "
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "2.174.11"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "111.85"
$ws.Range("E48").Value = "  -8.35%  "
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("E50").Value = "  +11.26%  "
$ws.Range("D51").Value = "0.944"
$ws.Range("E51").Value = "  -2.29%  "
